$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.385.06'
$ws.Range("E2").Value = '  -1.99%  '
$ws.Range("D3").Value = '1.795.96'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.24%  '
$ws.Range("E5").Value = '  -0.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '307.58'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.09%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4509'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.40%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3597'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.49%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.87'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.08%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07076'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.48%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8855'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.05%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07787'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.72%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.41'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.09%  '
$ws.Range("D14").Value = '1.824.43'
$ws.Range("E14").Value = '  +0.16%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.285'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.95%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.338'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.86%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '84.95'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.46%  '
$ws.Range("E18").Value = '  -0.14%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000008507'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.45%  '
$ws.Range("E20").Value = '  -0.13%  '
$ws.Range("E21").Value = '  -1.56%  '
$ws.Range("D22").Value = '26.401.07'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.991'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.23%  '
$ws.Range("B24").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C24").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D24").Value = '2.060.88'
$ws.Range("E24").Value = '  +0.76%  '
$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.53'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.90%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.973'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.40%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '152.13'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.66%  '
$ws.Range("E28").Value = '  -1.64%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.029'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.26%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '112.14'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.61%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.872'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.14%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08686'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.28%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.050'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.19%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.753'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +7.20%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.449'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.73%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7232'
$ws.Range("D36").Style = "Normal"
$ws.Range("E37").Value = '  -2.24%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.005'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.069'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.84%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01930'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.35%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.05097'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.864'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.28%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.892'
$ws.Range("D43").Style = "Normal"
$ws.Range("E44").Value = '  +1.64%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1511'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.52%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.016'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.53%  '
$ws.Range("E47").Value = '  -0.09%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4625'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '101.04'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.807'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.67%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.581'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.13%  '
